$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1: keep top+bottom border only (drop left & right) -> matches target borderId 4
$ws1.Range("C1").Borders.Item(7).LineStyle = -4142   # left
$ws1.Range("C1").Borders.Item(10).LineStyle = -4142  # right

# D1: keep top+bottom+right border only (drop left) -> matches target borderId 5
$ws1.Range("D1").Borders.Item(7).LineStyle = -4142   # left

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

# C1 / F1: keep top+bottom border only (drop left & right)
$ws2.Range("C1").Borders.Item(7).LineStyle = -4142
$ws2.Range("C1").Borders.Item(10).LineStyle = -4142
$ws2.Range("F1").Borders.Item(7).LineStyle = -4142
$ws2.Range("F1").Borders.Item(10).LineStyle = -4142

# D1 / G1: keep top+bottom+right border only (drop left)
$ws2.Range("D1").Borders.Item(7).LineStyle = -4142
$ws2.Range("G1").Borders.Item(7).LineStyle = -4142

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5 (model_size/change column had no value)
$ws2.Range("G5").ClearContents()
